$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-09 Saturday" "2024-11-10 Sunday"

Replace-Text "63×69=" "23×46="
Replace-Text "98×32=" "61×44="
Replace-Text "24×94=" "28×31="
Replace-Text "59×97=" "98×28="
Replace-Text "63×84=" "35×76="
Replace-Text "34×57=" "12×56="
Replace-Text "53×69=" "26×87="
Replace-Text "33×78=" "68×11="
Replace-Text "31×95=" "76×23="
Replace-Text "34×18=" "76×46="
Replace-Text "53×43=" "73×91="
Replace-Text "82×65=" "36×13="
Replace-Text "49×78=" "38×54="
Replace-Text "93×58=" "98×14="
Replace-Text "87×93=" "12×42="
Replace-Text "49×26=" "67×55="
Replace-Text "48×32=" "26×76="
Replace-Text "44×68=" "55×13="
Replace-Text "73×95=" "85×27="
Replace-Text "35×28=" "76×36="
Replace-Text "95×92=" "47×54="
Replace-Text "30×43=" "56×77="
Replace-Text "41×16=" "57×38="
Replace-Text "71×75=" "82×29="
Replace-Text "17×11=" "28×76="
